$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$rowCount = $usedRange.Rows.Count

for ($r = 1; $r -le $rowCount; $r++) {
    $weekCell = $ws.Cells.Item($r, 2)
    $weekVal = $weekCell.Value2
    if ($weekVal -ne $null -and $weekVal.ToString().StartsWith("Week 7")) {
        $dateCell = $ws.Cells.Item($r, 3)
        $dateVal = $dateCell.Value2
        if ($dateVal -eq $null -or $dateVal -eq "") {
            $dateCell.Value = "10/16"
        }
    }
}
